$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.437.95'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '3.082.67'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''521.54'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').Value = '''142.73'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''0.438'
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('E11').Value = '  +2.47%  '
$ws.Range('D12').Value = '3.610.32'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('D14').Value = '''26.61'
$ws.Range('E14').Value = '  +3.50%  '
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '58.461.40'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').Value = '3.080.96'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').Value = '''6.13'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('E19').Value = '  -1.84%  '
$ws.Range('D20').Value = '''8.10'
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('D21').Value = '''341.16'
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '''0.504'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').Value = '''65.69'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('E25').Value = '  -0.52%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').Value = '0.0₃0915'
$ws.Range('E27').Value = '  -1.64%  '
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('D29').Value = '''7.21'
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('E30').Value = '  +1.33%  '
$ws.Range('D31').Value = '''20.95'
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('E32').Value = '  +1.99%  '
$ws.Range('D33').Value = '''154.11'
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('D34').Value = '''4.59'
$ws.Range('E34').Value = '  +1.51%  '
$ws.Range('D35').Value = '''6.04'
$ws.Range('E35').Value = '  +2.24%  '
$ws.Range('D36').Value = '''26.90'
$ws.Range('E36').Value = '  -4.09%  '
$ws.Range('D37').Value = '''1.31'
$ws.Range('E37').Value = '  +5.67%  '
$ws.Range('D38').Value = '''0.0677'
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').Value = '3.124.21'
$ws.Range('E39').Value = '  -0.33%  '
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').Value = '''36.71'
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '''1.48'
$ws.Range('E42').Value = '  +7.20%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '''1.00'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '''0.665'
$ws.Range('E44').Value = '  -1.05%  '
$ws.Range('D45').Value = '2.267.45'
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('D46').Value = '''0.0256'
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('D47').Value = '''20.67'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('D48').Value = '''0.954'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('E50').Value = '  +7.28%  '
$ws.Range('D51').Value = '''264.37'
$ws.Range('E51').Value = '  +9.73%  '
